$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the row above so the new date cell keeps
# the same style (numFmt 14, "m/d/yyyy") instead of creating a new style.
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats

# New status entry row
$ws.Range("A21").Value = 43175
$ws.Range("B21").Value = "Android bugfix post methods -> research (20%)"
$ws.Range("C21").Value = "C# admin note communication between two windows fixed (80%)"
$ws.Range("D21").Value = "WebApp bugfix (100%)"

# Update selection to match new last row
$null = $ws.Range("D21").Select()
